$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FE_LFT_#1")
$ws.Cells.Item(62, 1).Value = 45848.49329861111
$ws.Cells.Item(62, 1).NumberFormat = $ws.Cells.Item(61, 1).NumberFormat
$ws.Cells.Item(62, 2).Value = "0x01,0x7c"
$ws.Cells.Item(62, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(62, 4).Value = "0x01,0x4C"
$ws.Cells.Item(62, 5).Value = "0xf"
$ws.Cells.Item(62, 6).Value = 380
$ws.Cells.Item(62, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(62, 8).Value = 332
$ws.Cells.Item(62, 9).Value = 15

$ws.Cells.Item(63, 1).Value = 45849.49347222222
$ws.Cells.Item(63, 1).NumberFormat = $ws.Cells.Item(62, 1).NumberFormat
$ws.Cells.Item(63, 2).Value = "0x01,0x7c"
$ws.Cells.Item(63, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(63, 4).Value = "0x01,0x4C"
$ws.Cells.Item(63, 5).Value = "0xf"
$ws.Cells.Item(63, 6).Value = 380
$ws.Cells.Item(63, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(63, 8).Value = 332
$ws.Cells.Item(63, 9).Value = 15

$ws.Cells.Item(64, 1).Value = 45850.49703703704
$ws.Cells.Item(64, 1).NumberFormat = $ws.Cells.Item(63, 1).NumberFormat
$ws.Cells.Item(64, 2).Value = "0x01,0x7c"
$ws.Cells.Item(64, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(64, 4).Value = "0x01,0x48"
$ws.Cells.Item(64, 5).Value = "0xf"
$ws.Cells.Item(64, 6).Value = 380
$ws.Cells.Item(64, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(64, 8).Value = 328
$ws.Cells.Item(64, 9).Value = 15


$ws = $wb.Worksheets.Item("FE_LFT_#2")
$ws.Cells.Item(62, 1).Value = 45848.49329861111
$ws.Cells.Item(62, 1).NumberFormat = $ws.Cells.Item(61, 1).NumberFormat
$ws.Cells.Item(62, 2).Value = "0x01,0x90"
$ws.Cells.Item(62, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(62, 4).Value = "0x01,0x5C"
$ws.Cells.Item(62, 5).Value = "0xe"
$ws.Cells.Item(62, 6).Value = 400
$ws.Cells.Item(62, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(62, 8).Value = 348
$ws.Cells.Item(62, 9).Value = 14

$ws.Cells.Item(63, 1).Value = 45849.49347222222
$ws.Cells.Item(63, 1).NumberFormat = $ws.Cells.Item(62, 1).NumberFormat
$ws.Cells.Item(63, 2).Value = "0x01,0x90"
$ws.Cells.Item(63, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(63, 4).Value = "0x01,0x5C"
$ws.Cells.Item(63, 5).Value = "0xe"
$ws.Cells.Item(63, 6).Value = 400
$ws.Cells.Item(63, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(63, 8).Value = 348
$ws.Cells.Item(63, 9).Value = 14

$ws.Cells.Item(64, 1).Value = 45850.49703703704
$ws.Cells.Item(64, 1).NumberFormat = $ws.Cells.Item(63, 1).NumberFormat
$ws.Cells.Item(64, 2).Value = "0x01,0x90"
$ws.Cells.Item(64, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(64, 4).Value = "0x01,0x58"
$ws.Cells.Item(64, 5).Value = "0xe"
$ws.Cells.Item(64, 6).Value = 400
$ws.Cells.Item(64, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(64, 8).Value = 344
$ws.Cells.Item(64, 9).Value = 14


$ws = $wb.Worksheets.Item("FE_PLT_#1")
$ws.Cells.Item(62, 1).Value = 45848.49329861111
$ws.Cells.Item(62, 1).NumberFormat = $ws.Cells.Item(61, 1).NumberFormat
$ws.Cells.Item(62, 2).Value = "0x00,0x6e"
$ws.Cells.Item(62, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(62, 4).Value = "0x00,0x66"
$ws.Cells.Item(62, 5).Value = "0x3"
$ws.Cells.Item(62, 6).Value = 110
$ws.Cells.Item(62, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(62, 8).Value = 102
$ws.Cells.Item(62, 9).Value = 3

$ws.Cells.Item(63, 1).Value = 45849.49347222222
$ws.Cells.Item(63, 1).NumberFormat = $ws.Cells.Item(62, 1).NumberFormat
$ws.Cells.Item(63, 2).Value = "0x00,0x6e"
$ws.Cells.Item(63, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(63, 4).Value = "0x00,0x66"
$ws.Cells.Item(63, 5).Value = "0x3"
$ws.Cells.Item(63, 6).Value = 110
$ws.Cells.Item(63, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(63, 8).Value = 102
$ws.Cells.Item(63, 9).Value = 3

$ws.Cells.Item(64, 1).Value = 45850.49703703704
$ws.Cells.Item(64, 1).NumberFormat = $ws.Cells.Item(63, 1).NumberFormat
$ws.Cells.Item(64, 2).Value = "0x00,0x6e"
$ws.Cells.Item(64, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(64, 4).Value = "0x00,0x65"
$ws.Cells.Item(64, 5).Value = "0x3"
$ws.Cells.Item(64, 6).Value = 110
$ws.Cells.Item(64, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(64, 8).Value = 101
$ws.Cells.Item(64, 9).Value = 3


$ws = $wb.Worksheets.Item("FE_PLT_#2")
$ws.Cells.Item(62, 1).Value = 45848.49329861111
$ws.Cells.Item(62, 1).NumberFormat = $ws.Cells.Item(61, 1).NumberFormat
$ws.Cells.Item(62, 2).Value = "0x00,0x6e"
$ws.Cells.Item(62, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(62, 4).Value = "0x00,0x66"
$ws.Cells.Item(62, 5).Value = "0x3"
$ws.Cells.Item(62, 6).Value = 110
$ws.Cells.Item(62, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(62, 8).Value = 102
$ws.Cells.Item(62, 9).Value = 3

$ws.Cells.Item(63, 1).Value = 45849.49347222222
$ws.Cells.Item(63, 1).NumberFormat = $ws.Cells.Item(62, 1).NumberFormat
$ws.Cells.Item(63, 2).Value = "0x00,0x6e"
$ws.Cells.Item(63, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(63, 4).Value = "0x00,0x66"
$ws.Cells.Item(63, 5).Value = "0x3"
$ws.Cells.Item(63, 6).Value = 110
$ws.Cells.Item(63, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(63, 8).Value = 102
$ws.Cells.Item(63, 9).Value = 3

$ws.Cells.Item(64, 1).Value = 45850.49703703704
$ws.Cells.Item(64, 1).NumberFormat = $ws.Cells.Item(63, 1).NumberFormat
$ws.Cells.Item(64, 2).Value = "0x00,0x6e"
$ws.Cells.Item(64, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(64, 4).Value = "0x00,0x65"
$ws.Cells.Item(64, 5).Value = "0x3"
$ws.Cells.Item(64, 6).Value = 110
$ws.Cells.Item(64, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(64, 8).Value = 101
$ws.Cells.Item(64, 9).Value = 3

